$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 8838.477000000001
$ws.Cells.Item(40, 9).Value = 7914.857
$ws.Cells.Item(40, 10).Value = 10685.714
$ws.Cells.Item(40, 11).Value = 7914.857
$ws.Cells.Item(40, 12).Value = 10685.714
$ws.Cells.Item(40, 13).Value = -7739.857
$ws.Cells.Item(40, 14).Value = -11035.714

$ws.Cells.Item(43, 8).Value = 2063.9092
$ws.Cells.Item(43, 9).Value = 1885.5714
$ws.Cells.Item(43, 10).Value = 2376
$ws.Cells.Item(43, 11).Value = 1885.5714
$ws.Cells.Item(43, 12).Value = 2376
$ws.Cells.Item(43, 13).Value = -1816.5714
$ws.Cells.Item(43, 14).Value = -2514

$ws.Cells.Item(75, 8).Value = 20000
$ws.Cells.Item(75, 10).Value = 20000
$ws.Cells.Item(75, 12).Value = 20000
$ws.Cells.Item(75, 14).Value = -21872

$ws.Cells.Item(78, 8).Value = 20000
$ws.Cells.Item(78, 10).Value = 20000
$ws.Cells.Item(78, 12).Value = 60000
$ws.Cells.Item(78, 14).Value = -69360

$ws.Cells.Item(112, 8).Value = 2448.394
$ws.Cells.Item(112, 10).Value = 2448.394
$ws.Cells.Item(112, 12).Value = 7345.181999999999
$ws.Cells.Item(112, 14).Value = -9561.181999999999

$ws.Cells.Item(137, 8).Value = 2776.362
$ws.Cells.Item(137, 10).Value = 4816.0835
$ws.Cells.Item(137, 12).Value = 14448.2505
$ws.Cells.Item(137, 14).Value = -19548.2505

$ws.Cells.Item(138, 8).Value = 6057.5933
$ws.Cells.Item(138, 10).Value = 7352.0654
$ws.Cells.Item(138, 12).Value = 22056.1962
$ws.Cells.Item(138, 14).Value = -32336.1962

$ws.Cells.Item(141, 8).Value = 5281.273
$ws.Cells.Item(141, 9).Value = 5249.4375
$ws.Cells.Item(141, 11).Value = 15748.3125
$ws.Cells.Item(141, 13).Value = -10568.3125

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3223.3977
$ws.Cells.Item(32, 9).Value = 3223.3977
$ws.Cells.Item(32, 11).Value = 3223.3977
$ws.Cells.Item(32, 13).Value = -2936.3977

$ws.Cells.Item(45, 8).Value = 2525.389
$ws.Cells.Item(45, 9).Value = 2075.4546
$ws.Cells.Item(45, 11).Value = 2075.4546
$ws.Cells.Item(45, 13).Value = -1698.4546

$ws.Cells.Item(74, 8).Value = 1759.4032
$ws.Cells.Item(74, 9).Value = 1779.6227
$ws.Cells.Item(74, 11).Value = 1779.6227
$ws.Cells.Item(74, 13).Value = -905.6226999999999

$ws.Cells.Item(77, 8).Value = 1759.4032
$ws.Cells.Item(77, 9).Value = 1779.6227
$ws.Cells.Item(77, 11).Value = 8898.113499999999
$ws.Cells.Item(77, 13).Value = -4530.113499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 18377.61
$ws.Cells.Item(134, 9).Value = 2663.984
$ws.Cells.Item(134, 10).Value = 505500
$ws.Cells.Item(134, 11).Value = 7991.951999999999
$ws.Cells.Item(134, 12).Value = 1516500
$ws.Cells.Item(134, 13).Value = -5456.951999999999
$ws.Cells.Item(134, 14).Value = -1521570

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 686.3333
$ws.Cells.Item(22, 9).Value = 428.45456
$ws.Cells.Item(22, 11).Value = 428.45456
$ws.Cells.Item(22, 13).Value = -78.45456000000001

$ws.Cells.Item(132, 8).Value = 3793.36
$ws.Cells.Item(132, 9).Value = 2202
$ws.Cells.Item(132, 10).Value = 7175
$ws.Cells.Item(132, 11).Value = 6606
$ws.Cells.Item(132, 12).Value = 21525
$ws.Cells.Item(132, 13).Value = -4076
$ws.Cells.Item(132, 14).Value = -26585

$ws.Cells.Item(134, 8).Value = 305140.47
$ws.Cells.Item(134, 9).Value = 2283.9656
$ws.Cells.Item(134, 11).Value = 6851.8968
$ws.Cells.Item(134, 13).Value = -4316.8968

$ws.Cells.Item(135, 8).Value = 49997.184
$ws.Cells.Item(135, 10).Value = 49997.184
$ws.Cells.Item(135, 12).Value = 49997.184
$ws.Cells.Item(135, 14).Value = -60137.184

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(64, 8).Value = 166667260
$ws.Cells.Item(64, 9).Value = 166667260
$ws.Cells.Item(64, 10).Value = 0
$ws.Cells.Item(64, 11).Value = 500001780
$ws.Cells.Item(64, 12).Value = 0
$ws.Cells.Item(64, 13).Value = -500001510
$ws.Cells.Item(64, 14).ClearContents()

$ws.Cells.Item(67, 8).Value = 166667260
$ws.Cells.Item(67, 9).Value = 166667260
$ws.Cells.Item(67, 10).Value = 0
$ws.Cells.Item(67, 11).Value = 500001780
$ws.Cells.Item(67, 12).Value = 0
$ws.Cells.Item(67, 13).Value = -500000844
$ws.Cells.Item(67, 14).ClearContents()

$ws.Cells.Item(120, 8).Value = 257494.5
$ws.Cells.Item(120, 9).Value = 257494.5
$ws.Cells.Item(120, 10).Value = 0
$ws.Cells.Item(120, 11).Value = 772483.5
$ws.Cells.Item(120, 12).Value = 0
$ws.Cells.Item(120, 13).Value = -767645.5
$ws.Cells.Item(120, 14).ClearContents()

$ws.Cells.Item(124, 8).Value = 251572
$ws.Cells.Item(124, 9).Value = 251572
$ws.Cells.Item(124, 11).Value = 754716
$ws.Cells.Item(124, 13).Value = -749806

$ws.Cells.Item(129, 8).Value = 68370.8
$ws.Cells.Item(129, 9).Value = 751.8333
$ws.Cells.Item(129, 10).Value = 113450.11
$ws.Cells.Item(129, 11).Value = 2255.4999
$ws.Cells.Item(129, 12).Value = 340350.33
$ws.Cells.Item(129, 13).Value = 2744.5001
$ws.Cells.Item(129, 14).Value = -350350.33

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(18, 8).Value = 1000000000
$ws.Cells.Item(18, 10).Value = 1000000000
$ws.Cells.Item(18, 12).Value = 1000000000
$ws.Cells.Item(18, 14).Value = -1000000586

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 171167.67
$ws.Cells.Item(40, 9).Value = 252751.5
$ws.Cells.Item(40, 11).Value = 252751.5
$ws.Cells.Item(40, 13).Value = -252615.5

$ws.Cells.Item(46, 8).Value = 2465.0557
$ws.Cells.Item(46, 9).Value = 2592.2307
$ws.Cells.Item(46, 10).Value = 2134.4
$ws.Cells.Item(46, 11).Value = 2592.2307
$ws.Cells.Item(46, 12).Value = 2134.4
$ws.Cells.Item(46, 13).Value = -2404.2307
$ws.Cells.Item(46, 14).Value = -2510.4

$ws.Cells.Item(68, 8).Value = 251724.75
$ws.Cells.Item(68, 10).Value = 334966.66
$ws.Cells.Item(68, 12).Value = 334966.66
$ws.Cells.Item(68, 14).Value = -336464.66

$ws.Cells.Item(71, 8).Value = 251724.75
$ws.Cells.Item(71, 10).Value = 334966.66
$ws.Cells.Item(71, 12).Value = 1674833.3
$ws.Cells.Item(71, 14).Value = -1682321.3

$ws.Cells.Item(132, 8).Value = 4503.175
$ws.Cells.Item(132, 9).Value = 3982.6365
$ws.Cells.Item(132, 11).Value = 11947.9095
$ws.Cells.Item(132, 13).Value = -9417.9095

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 77762
$ws.Cells.Item(62, 9).Value = 204880.4
$ws.Cells.Item(62, 11).Value = 204880.4
$ws.Cells.Item(62, 13).Value = -204256.4

$ws.Cells.Item(65, 8).Value = 77762
$ws.Cells.Item(65, 9).Value = 204880.4
$ws.Cells.Item(65, 11).Value = 1024402
$ws.Cells.Item(65, 13).Value = -1021282

$ws.Cells.Item(80, 8).Value = 206803.25
$ws.Cells.Item(80, 10).Value = 206803.25
$ws.Cells.Item(80, 12).Value = 206803.25
$ws.Cells.Item(80, 14).Value = -208799.25

$ws.Cells.Item(83, 8).Value = 206803.25
$ws.Cells.Item(83, 10).Value = 206803.25
$ws.Cells.Item(83, 12).Value = 620409.75
$ws.Cells.Item(83, 14).Value = -630393.75

$ws.Cells.Item(126, 8).Value = 1675.4762
$ws.Cells.Item(126, 9).Value = 1681.6923
$ws.Cells.Item(126, 11).Value = 5045.0769
$ws.Cells.Item(126, 13).Value = -2575.0769

$ws.Cells.Item(132, 8).Value = 34701.766
$ws.Cells.Item(132, 9).Value = 1497.909
$ws.Cells.Item(132, 11).Value = 4493.727000000001
$ws.Cells.Item(132, 13).Value = -1963.727000000001

$ws.Cells.Item(136, 8).Value = 282233.3
$ws.Cells.Item(136, 9).Value = 289003.38
$ws.Cells.Item(136, 11).Value = 867010.14
$ws.Cells.Item(136, 13).Value = -864460.14
